# "Filter to query only 2023 data included"
#
# 1. Negate every non-blank value in column Q (interest_cover) for data rows
#    2-61, and turn the one blank Q cell (row 41) into an explicit 0.
# 2. Remove the stray summary row 63 (D63 = MAX(D2:D61)) that was left
#    below the data table - this also shrinks the sheet's used range/
#    dimension back down to A1:U61.
# 3. Reset the view back to the top-left corner (A1) so there is no
#    lingering scrolled/selected state pointing past the trimmed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Flip the sign of the interest_cover column (Q) for rows 2-61 ---
for ($r = 2; $r -le 61; $r++) {
    $cell = $ws.Cells.Item($r, 17)   # column Q
    $v = $cell.Value2
    if ($v -eq $null) {
        $cell.Value = 0
    } else {
        $cell.Value = -$v
    }
}

# --- 2. Delete the extra MAX() summary row below the table ---
$ws.Rows.Item(63).Delete()

# --- 3. Return the view to A1 (clears the stale scroll/selection state) ---
$ws.Range("A1").Select()
